$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Prime formatting first:
#    - Copy the existing plain data style (from B2, which already uses
#      the non-bold Arial/theme-colour style) onto every cell that will
#      receive new content, so new cells line up with the sheet's
#      existing "data row" look instead of falling back to a generic
#      default style.
# ------------------------------------------------------------------
$ws.Range("B2").Copy($ws.Range("D2"))
$ws.Range("B2").Copy($ws.Range("A3:D11"))

# Apply the date number format (m/d/yyyy) to the Assigned date column
# and to the two Complete Date cells that hold real dates.
$ws.Range("C3:C11").NumberFormat = "m/d/yyyy"
$ws.Range("D3:D4").NumberFormat = "m/d/yyyy"

# Column widths (B and C got a little narrower in the source edit)
$ws.Range("B1").ColumnWidth = 29.3
$ws.Range("C1").ColumnWidth = 11.55

# ------------------------------------------------------------------
# 2) Header row - add the new "In progress" status header in D2
# ------------------------------------------------------------------
$ws.Range("D2").Value = "In progress"

# ------------------------------------------------------------------
# 3) Data rows
# ------------------------------------------------------------------
# Row 3: Arpit / Make login page (completed)
$ws.Range("A3").Value = "Arpit "
$ws.Range("B3").Value = "Make login page"
$ws.Range("C3").Value = 44594
$ws.Range("D3").Value = 44645

# Row 4: Arpit / Make registration page (completed)
$ws.Range("A4").Value = "Arpit "
$ws.Range("B4").Value = "Make registration page"
$ws.Range("C4").Value = 44594
$ws.Range("D4").Value = 44645

# Row 5: Ishika / Make forgot password page (in progress)
$ws.Range("A5").Value = "Ishika"
$ws.Range("B5").Value = "Make forgot password page"
$ws.Range("C5").Value = 44594
$ws.Range("D5").Value = "In progress"

# Row 6: Ishika / Design report on billing (in progress)
$ws.Range("A6").Value = "Ishika"
$ws.Range("B6").Value = "Design report on billing"
$ws.Range("C6").Value = 44638
$ws.Range("D6").Value = "In progress"

# Row 7: Maeve / Make upcoming trips (in progress)
$ws.Range("A7").Value = "Maeve"
$ws.Range("B7").Value = "Make upcoming trips "
$ws.Range("C7").Value = 44638
$ws.Range("D7").Value = "In progress"

# Row 8: Maeve / Edit & Delete reservations (in progress, merged task)
$ws.Range("A8").Value = "Maeve"
$ws.Range("B8").Value = "Edit & Delete reservations"
$ws.Range("C8").Value = 44638
$ws.Range("D8").Value = "In progress"

# Row 9: Tri / Account Page (new task, in progress)
$ws.Range("A9").Value = "Tri"
$ws.Range("B9").Value = "Account Page"
$ws.Range("C9").Value = 44639
$ws.Range("D9").Value = "In progress"

# Row 10: Tri / Test the pages (in progress)
$ws.Range("A10").Value = "Tri"
$ws.Range("B10").Value = "Test the pages"
$ws.Range("C10").Value = 44594
$ws.Range("D10").Value = "In progress"

# Row 11 (new row): Arpit / Search Page (in progress)
$ws.Range("A11").Value = "Arpit "
$ws.Range("B11").Value = "Search Page"
$ws.Range("C11").Value = 44638
$ws.Range("D11").Value = "In progress"
